$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EE")

# New E-column (5th column) values for rows 3..63, keyed by row number.
# Rows 46 and 62 are intentionally absent (their E value is unchanged).
$eValues = @{
    3  = 15
    4  = 8
    5  = 20
    6  = 16
    7  = 15
    8  = 17
    9  = 5
    10 = 11
    11 = 10
    12 = 11
    13 = 11
    14 = 11
    15 = 20
    16 = 16
    17 = 11
    18 = 14
    19 = 18
    20 = 16
    21 = 12
    22 = 18
    23 = 8
    24 = 17
    25 = 8
    26 = 8
    27 = 7
    28 = 18
    29 = 20
    30 = 20
    31 = 17
    32 = 7
    33 = 13
    34 = 5
    35 = 12
    36 = 18
    37 = 6
    38 = 5
    39 = 5
    40 = 7
    41 = 10
    42 = 9
    43 = 17
    44 = 20
    45 = 9
    47 = 14
    48 = 9
    49 = 7
    50 = 7
    51 = 10
    52 = 14
    53 = 17
    54 = 13
    55 = 20
    56 = 19
    57 = 6
    58 = 13
    59 = 5
    60 = 13
    61 = 8
    63 = 11
}

for ($row = 3; $row -le 63; $row++) {
    $aCell = $ws.Cells.Item($row, 1)
    if ($aCell.Value2 -ne $null) {
        $aCell.Value2 = $aCell.Value2 + 20000
    }

    if ($eValues.ContainsKey($row)) {
        $ws.Cells.Item($row, 5).Value2 = $eValues[$row]
    }
}
